$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from A12 onto A13 and A14 so they share the
# same numFmt (ddd dd/mm/yyyy) style index as the other date cells.
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)

# Row 13: Date 41440 (2013-06-15), Effort 2.5, Task "Implementation tc12"
$ws.Range("A13").Value = 41440
$ws.Range("B13").Value = 2.5
$ws.Range("D13").Value = "Implementation tc12"

# Row 14: Date 41442 (2013-06-17), Effort 2, Task "Implementation tc12"
$ws.Range("A14").Value = 41442
$ws.Range("B14").Value = 2
$ws.Range("D14").Value = "Implementation tc12"

# Update selection to match diff (active cell A14)
$ws.Range("A14").Select()
